$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column H - copy formatting from the neighboring header (G1)
# so it matches the bold/centered/bordered header style, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Label"

# Label values: 0 for Control rows, 1 for MDD rows (two repeated blocks of 10 rows each)
$labels = @(0, 0, 0, 0, 0, 1, 1, 1, 1, 1, 0, 0, 0, 0, 0, 1, 1, 1, 1, 1)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $labels[$i]
}
